$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 08:22"

# --- Corea del Sur / Austria-ish row (row 19) new totals ---
$ws.Range("B19").Value = 13271
$ws.Range("C19").Value = 27
$ws.Range("E19").Value = 7736

# --- Lituania (row 68) new totals ---
$ws.Range("E68").Value = 974
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 17

# --- Countries resorted: Afganistan overtakes Cuba (rows 88/89) ---
# Row 88 now shows Afganistan with its updated totals
$ws.Range("A88").Value = "Afganistan"
$ws.Range("B88").Value = 521
$ws.Range("C88").Value = 37
$ws.Range("D88").Value = 32
$ws.Range("E88").Value = 474
$ws.Range("F88").Value = 0

# Row 89 now shows Cuba (its totals are unchanged from before the resort)
$ws.Range("A89").Value = "Cuba"
$ws.Range("B89").Value = 515
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 28
$ws.Range("E89").Value = 472
$ws.Range("F89").Value = 15

# --- Taiwan (row 97) new totals ---
$ws.Range("B97").Value = 382
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 91
$ws.Range("E97").Value = 285
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 6

# --- Vietnam (row 109) new totals ---
$ws.Range("D109").Value = 144
$ws.Range("E109").Value = 111

# --- Countries resorted: Mayotte overtakes Kenia (rows 117/118) ---
# Row 117 now shows Mayotte with its updated totals
$ws.Range("A117").Value = "Mayotte"
$ws.Range("B117").Value = 186
$ws.Range("C117").Value = 2
$ws.Range("D117").Value = 26
$ws.Range("E117").Value = 158
$ws.Range("F117").Value = 4
$ws.Range("H117").Value = 2

# Row 118 now shows Kenia (its totals are unchanged from before the resort)
$ws.Range("A118").Value = "Kenia"
$ws.Range("D118").Value = 12
$ws.Range("E118").Value = 165
$ws.Range("F118").Value = 2
$ws.Range("H118").Value = 7
